# Auto-generated edit script applying Durandal_Profits.xlsx market-price refresh
# updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N) per sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 2804.96
$ws.Range("I113").Value = 2867.2942
$ws.Range("J113").Value = 2672.5
$ws.Range("K113").Value = 2867.2942
$ws.Range("L113").Value = 2672.5
$ws.Range("M113").Value = 386.7058000000002
$ws.Range("N113").Value = -9180.5
$ws.Range("H132").Value = 5209765
$ws.Range("I132").Value = 5815458.5
$ws.Range("J132").Value = 799.8
$ws.Range("K132").Value = 17446375.5
$ws.Range("L132").Value = 2399.4
$ws.Range("M132").Value = -17443845.5
$ws.Range("N132").Value = -7459.4
$ws.Range("H137").Value = 1171.826
$ws.Range("I137").Value = 1133.9333
$ws.Range("J137").Value = 1242.875
$ws.Range("K137").Value = 3401.7999
$ws.Range("L137").Value = 3728.625
$ws.Range("M137").Value = -851.7999
$ws.Range("N137").Value = -8828.625
$ws.Range("H138").Value = 5347.5356
$ws.Range("I138").Value = 4719.4
$ws.Range("J138").Value = 5484.087
$ws.Range("K138").Value = 14158.2
$ws.Range("L138").Value = 16452.261
$ws.Range("M138").Value = -9018.199999999999
$ws.Range("N138").Value = -26732.261

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 415767.44
$ws.Range("I32").Value = 4058.4746
$ws.Range("J32").Value = 2624024.5
$ws.Range("K32").Value = 4058.4746
$ws.Range("L32").Value = 2624024.5
$ws.Range("M32").Value = -3771.4746
$ws.Range("N32").Value = -2624598.5
$ws.Range("H45").Value = 2368.7334
$ws.Range("I45").Value = 2441.375
$ws.Range("J45").Value = 2285.7144
$ws.Range("K45").Value = 2441.375
$ws.Range("L45").Value = 2285.7144
$ws.Range("M45").Value = -2064.375
$ws.Range("N45").Value = -3039.7144

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2911.7778
$ws.Range("I20").Value = 2170.5293
$ws.Range("J20").Value = 4171.9
$ws.Range("K20").Value = 2170.5293
$ws.Range("L20").Value = 4171.9
$ws.Range("M20").Value = -1923.5293
$ws.Range("N20").Value = -4665.9
$ws.Range("H22").Value = 315.66666
$ws.Range("I22").Value = 320.125
$ws.Range("J22").Value = 280
$ws.Range("K22").Value = 320.125
$ws.Range("L22").Value = 280
$ws.Range("M22").Value = -147.125
$ws.Range("N22").Value = -626
$ws.Range("H140").Value = 51193.332
$ws.Range("J140").Value = 51193.332
$ws.Range("L140").Value = 51193.332
$ws.Range("N140").Value = -61553.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H38").Value = 50035020
$ws.Range("I38").Value = 50035020
$ws.Range("K38").Value = 50035020
$ws.Range("M38").Value = -50034643
$ws.Range("H46").Value = 50035020
$ws.Range("I46").Value = 50035020
$ws.Range("K46").Value = 50035020
$ws.Range("M46").Value = -50034809

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 883.2059
$ws.Range("I5").Value = 411.94446
$ws.Range("J5").Value = 1413.375
$ws.Range("K5").Value = 1235.83338
$ws.Range("L5").Value = 4240.125
$ws.Range("M5").Value = -1123.83338
$ws.Range("N5").Value = -4464.125
$ws.Range("H68").Value = 1121.7255
$ws.Range("I68").Value = 780.1111
$ws.Range("J68").Value = 1308.0605
$ws.Range("K68").Value = 2340.3333
$ws.Range("L68").Value = 3924.1815
$ws.Range("M68").Value = -1529.3333
$ws.Range("N68").Value = -5546.181500000001
$ws.Range("H70").Value = 912
$ws.Range("I70").Value = 912
$ws.Range("K70").Value = 2736
$ws.Range("M70").Value = -2421
$ws.Range("H71").Value = 1121.7255
$ws.Range("I71").Value = 780.1111
$ws.Range("J71").Value = 1308.0605
$ws.Range("K71").Value = 7020.9999
$ws.Range("L71").Value = 11772.5445
$ws.Range("M71").Value = -2964.9999
$ws.Range("N71").Value = -19884.5445
$ws.Range("H73").Value = 912
$ws.Range("I73").Value = 912
$ws.Range("K73").Value = 2736
$ws.Range("M73").Value = -1644
$ws.Range("H107").Value = 1294.6207
$ws.Range("I107").Value = 664.6
$ws.Range("J107").Value = 1626.2106
$ws.Range("K107").Value = 1993.8
$ws.Range("L107").Value = 4878.6318
$ws.Range("M107").Value = -73.80000000000018
$ws.Range("N107").Value = -8718.631799999999
$ws.Range("H113").Value = 1009.68054
$ws.Range("I113").Value = 851.3333
$ws.Range("J113").Value = 1024.0758
$ws.Range("K113").Value = 2553.9999
$ws.Range("L113").Value = 3072.2274
$ws.Range("M113").Value = -383.9998999999998
$ws.Range("N113").Value = -7412.2274
$ws.Range("H135").Value = 883.2059
$ws.Range("I135").Value = 411.94446
$ws.Range("J135").Value = 1413.375
$ws.Range("K135").Value = 3707.50014
$ws.Range("L135").Value = 12720.375
$ws.Range("M135").Value = -1172.50014
$ws.Range("N135").Value = -17790.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1367.125
$ws.Range("I113").Value = 1057.4445
$ws.Range("J113").Value = 1765.2858
$ws.Range("K113").Value = 1057.4445
$ws.Range("L113").Value = 1765.2858
$ws.Range("M113").Value = 1112.5555
$ws.Range("N113").Value = -6105.2858
$ws.Range("H132").Value = 26617
$ws.Range("I132").Value = 2150.1562
$ws.Range("J132").Value = 113610.22
$ws.Range("K132").Value = 6450.4686
$ws.Range("L132").Value = 340830.66
$ws.Range("M132").Value = -3920.4686
$ws.Range("N132").Value = -345890.66

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 937.9231
$ws.Range("I16").Value = 803.3
$ws.Range("K16").Value = 803.3
$ws.Range("M16").Value = -633.3
$ws.Range("H30").Value = 500
$ws.Range("I30").Value = 500
$ws.Range("K30").Value = 500
$ws.Range("M30").Value = -392
$ws.Range("H136").Value = 4074.08
$ws.Range("I136").Value = 5297.44
$ws.Range("K136").Value = 15892.32
$ws.Range("M136").Value = -13342.32

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4595.625
$ws.Range("I81").Value = 4104.2856
$ws.Range("J81").Value = 4977.778
$ws.Range("K81").Value = 8208.5712
$ws.Range("L81").Value = 9955.556
$ws.Range("M81").Value = -7147.5712
$ws.Range("N81").Value = -12077.556
$ws.Range("H84").Value = 4595.625
$ws.Range("I84").Value = 4104.2856
$ws.Range("J84").Value = 4977.778
$ws.Range("K84").Value = 41042.856
$ws.Range("L84").Value = 49777.78
$ws.Range("M84").Value = -35738.856
$ws.Range("H136").Value = 35267.035
$ws.Range("I136").Value = 63179.625
$ws.Range("K136").Value = 189538.875
$ws.Range("M136").Value = -186988.875

